# Adds a new column (hours_to_drop) and a new location row (Cave of Shadows)
# to the Locations data-import worksheet, along with its associated quest
# item / quest reward text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column O: hours_to_drop -------------------------------------------------
$ws.Cells.Item(1, 15).Value = "hours_to_drop"

# Existing rows 2-70 default to 0 hours_to_drop.
for ($r = 2; $r -le 70; $r++) {
    $ws.Cells.Item($r, 15).Value = 0
}

# --- New row 71: Cave of Shadows --------------------------------------------------
$ws.Cells.Item(71, 1).Value = "Cave of Shadows"
$ws.Cells.Item(71, 2).Value = "Twisted Memories"
$ws.Cells.Item(71, 4).Value = "Ever Burning Candle"
$ws.Cells.Item(71, 5).Value = "A cave full of shadows of the deep, shadows of the mind, shadows of what use to be, what never was, what could have been ..."
$ws.Cells.Item(71, 9).Value = 1
$ws.Cells.Item(71, 10).Value = 1408
$ws.Cells.Item(71, 11).Value = 640
$ws.Cells.Item(71, 12).Value = 11
$ws.Cells.Item(71, 13).Value = "No"
$ws.Cells.Item(71, 15).Value = 1
